$null = 1
